# Generate Report for handback
# Updates the zh-cn and de-de status sheets: flips "Ready for handoff" rows to
# "Handed back: in sync with en-US", fills in the "Latest Target File" /
# "Latest Handback File" columns (E/F) with hyperlinked file names, and stamps
# the "Latest Handback DateTime" column (G) with the handback timestamp.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = $statusText
$ws.Range("B3").Value = $statusText

$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/cf085a11dc5af2db7c26a8aa4e891c792fd429cd/e2e/025500e4-db87-4bae-9d38-ba4f7073d4f3.md", "", "", "025500e4-db87-4bae-9d38-ba4f7073d4f3.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/07b48f26ebca8999451d8a5bd35cdc39a2b610c8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/025500e4-db87-4bae-9d38-ba4f7073d4f3.5d1d3d67c648da9ceb0c050f07fd60a45c2a5d64.zh-cn.xlf", "", "", "025500e4-db87-4bae-9d38-ba4f7073d4f3.5d1d3d67c648da9ceb0c050f07fd60a45c2a5d64.zh-cn.xlf")
$ws.Range("G2").Value = "2016-01-26 06:17:11"

$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/cf085a11dc5af2db7c26a8aa4e891c792fd429cd/e2e/772dba4e-725d-4d34-9189-60a88700dc9d.md", "", "", "772dba4e-725d-4d34-9189-60a88700dc9d.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/07b48f26ebca8999451d8a5bd35cdc39a2b610c8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/772dba4e-725d-4d34-9189-60a88700dc9d.ec57e943b01be122b7e40be922a375bf124c07e1.zh-cn.xlf", "", "", "772dba4e-725d-4d34-9189-60a88700dc9d.ec57e943b01be122b7e40be922a375bf124c07e1.zh-cn.xlf")
$ws.Range("G3").Value = "2016-01-26 06:17:11"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = $statusText
$ws.Range("B3").Value = $statusText

$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/cf085a11dc5af2db7c26a8aa4e891c792fd429cd/e2e/025500e4-db87-4bae-9d38-ba4f7073d4f3.md", "", "", "025500e4-db87-4bae-9d38-ba4f7073d4f3.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/151d4201d7eb4efe27845b288721af3ff209966f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/025500e4-db87-4bae-9d38-ba4f7073d4f3.5d1d3d67c648da9ceb0c050f07fd60a45c2a5d64.de-de.xlf", "", "", "025500e4-db87-4bae-9d38-ba4f7073d4f3.5d1d3d67c648da9ceb0c050f07fd60a45c2a5d64.de-de.xlf")
$ws.Range("G2").Value = "2016-01-26 06:17:33"

$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/cf085a11dc5af2db7c26a8aa4e891c792fd429cd/e2e/772dba4e-725d-4d34-9189-60a88700dc9d.md", "", "", "772dba4e-725d-4d34-9189-60a88700dc9d.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/151d4201d7eb4efe27845b288721af3ff209966f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/772dba4e-725d-4d34-9189-60a88700dc9d.ec57e943b01be122b7e40be922a375bf124c07e1.de-de.xlf", "", "", "772dba4e-725d-4d34-9189-60a88700dc9d.ec57e943b01be122b7e40be922a375bf124c07e1.de-de.xlf")
$ws.Range("G3").Value = "2016-01-26 06:17:33"
